# Update gh-pages output data: bump "想去人数" (want-to-go count) values
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 188
$ws1.Range("F8").Value = 267
$ws1.Range("F15").Value = 13004
$ws1.Range("F18").Value = 5347

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 188
$ws4.Range("F9").Value = 267
$ws4.Range("F17").Value = 13004
$ws4.Range("F21").Value = 5347
